$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Range("F1:J20").Interior.ColorIndex = -4142
Write-Host "done"
